# Weekly price update: insert a new most-recent-week record for
# "Feria Lagunitas de Puerto Montt" / Coliflor, pushing the existing
# history (rows 483-503) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 483; Excel shifts rows 483-503 down to 484-504
# (cell formatting, e.g. the date style on column D, travels with them).
$ws.Rows(483).Insert()

# Populate the newly inserted row 483 with this week's data.
$ws.Range("A483").Value = 4
$ws.Range("B483").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C483").Value = "Los Lagos"
$ws.Range("D483").Value = 45008
$ws.Range("E483").Value = 10
$ws.Range("F483").Value = 100112008
$ws.Range("G483").Value = "Coliflor"
$ws.Range("H483").Value = "Sin especificar"
$ws.Range("I483").Value = "Primera"
$ws.Range("J483").Value = 500
$ws.Range("K483").Value = 1700
$ws.Range("L483").Value = 1700
$ws.Range("M483").Value = 1700
$ws.Range("N483").Value = "$/unidad"
$ws.Range("O483").Value = "Región Metropolitana"
$ws.Range("P483").Value = 1700
$ws.Range("Q483").Value = 1
$ws.Range("R483").Value = "Hortaliza"
